# Auto-generated edit script: updates the cryptos price/volume table
# to match the Sun Jul 28 21:44:23 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.047.36"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "3.263.03"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'583.39"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'184.69"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.601"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("D11").Value = "'0.408"
$ws.Range("E11").Value = "  -3.71%  "
$ws.Range("D12").Value = "3.830.42"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").Value = "'27.39"
$ws.Range("E14").Value = "  -4.17%  "
$ws.Range("D15").Value = "68.034.95"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").Value = "'0.0000168"
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("D17").Value = "3.266.31"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "'5.71"
$ws.Range("E18").Value = "  -2.65%  "
$ws.Range("D19").Value = "'13.27"
$ws.Range("E19").Value = "  -2.85%  "
$ws.Range("D20").Value = "'417.68"
$ws.Range("E20").Value = "  +5.53%  "
$ws.Range("D21").Value = "'7.53"
$ws.Range("E21").Value = "  -2.79%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'71.15"
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("D24").Value = "'0.508"
$ws.Range("E24").Value = "  -2.76%  "
$ws.Range("E25").Value = "  -3.60%  "
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("D27").Value = "'9.37"
$ws.Range("E27").Value = "  -3.83%  "
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").Value = "'1.95"
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("D30").Value = "'22.56"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D31").Value = "'5.45"
$ws.Range("E31").Value = "  -5.07%  "
$ws.Range("D32").Value = "'6.85"
$ws.Range("E32").Value = "  -4.78%  "
$ws.Range("E33").Value = "  -5.08%  "
$ws.Range("D34").Value = "'162.82"
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("D35").Value = "'1.44"
$ws.Range("E35").Value = "  -5.47%  "
$ws.Range("D36").Value = "'1.88"
$ws.Range("E36").Value = "  -4.94%  "
$ws.Range("D37").Value = "'26.83"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "'0.795"
$ws.Range("E38").Value = "  -4.21%  "
$ws.Range("D39").Value = "'4.44"
$ws.Range("E39").Value = "  -4.01%  "
$ws.Range("E40").Value = "  -5.41%  "
$ws.Range("D41").Value = "2.633.59"
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("D42").Value = "'2.43"
$ws.Range("E42").Value = "  -5.30%  "
$ws.Range("D43").Value = "'0.0673"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("D44").Value = "'336.92"
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("D45").Value = "'24.22"
$ws.Range("E45").Value = "  -5.73%  "
$ws.Range("D46").Value = "'0.0273"
$ws.Range("E46").Value = "  -3.82%  "
$ws.Range("D47").Value = "'6.23"
$ws.Range("E47").Value = "  -2.37%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'0.975"
$ws.Range("E48").Value = "  -2.60%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.100"
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").Value = "'30.49"
$ws.Range("E51").Value = "  -5.53%  "
